# Auto-generated data refresh for Kujata_Profits workbook
# Updates currentAveragePrice / Leve profit columns (H:N) per sheet,
# matching the scheduled-runner price-refresh commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 88.75
$ws.Range("I9").Value = 65
$ws.Range("K9").Value = 65
$ws.Range("M9").Value = 104

$ws.Range("H15").Value = 2764.7954
$ws.Range("I15").Value = 2764.7954
$ws.Range("K15").Value = 8294.386200000001
$ws.Range("M15").Value = -8125.386200000001

$ws.Range("H100").Value = 2283.5
$ws.Range("I100").Value = 2161.3333
$ws.Range("J100").Value = 2650
$ws.Range("K100").Value = 2161.3333
$ws.Range("L100").Value = 2650
$ws.Range("M100").Value = -1620.3333
$ws.Range("N100").Value = -3732

$ws.Range("H138").Value = 454642.9
$ws.Range("I138").Value = 1180.2778
$ws.Range("J138").Value = 559288.1
$ws.Range("K138").Value = 3540.8334
$ws.Range("L138").Value = 1677864.3
$ws.Range("M138").Value = 1599.1666
$ws.Range("N138").Value = -1688144.3

$ws.Range("H141").Value = 2566.125
$ws.Range("I141").Value = 2504.1428
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 7512.428400000001
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -2332.428400000001
$ws.Range("N141").Value = -19360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 33345326
$ws.Range("I33").Value = 33345326
$ws.Range("K33").Value = 33345326
$ws.Range("M33").Value = -33344997

$ws.Range("H61").Value = 52632844
$ws.Range("I61").Value = 71429540
$ws.Range("J61").Value = 2116
$ws.Range("K61").Value = 71429540
$ws.Range("L61").Value = 2116
$ws.Range("M61").Value = -71429328
$ws.Range("N61").Value = -2540

$ws.Range("H122").Value = 1506.8379
$ws.Range("I122").Value = 1436.6666
$ws.Range("J122").Value = 1636.3846
$ws.Range("K122").Value = 4309.9998
$ws.Range("L122").Value = 4909.1538
$ws.Range("M122").Value = -1859.9998
$ws.Range("N122").Value = -9809.1538

$ws.Range("H132").Value = 1616.3954
$ws.Range("I132").Value = 1216.7576
$ws.Range("J132").Value = 2935.2
$ws.Range("K132").Value = 3650.2728
$ws.Range("L132").Value = 8805.599999999999
$ws.Range("M132").Value = -1120.2728
$ws.Range("N132").Value = -13865.6

$ws.Range("H136").Value = 52632844
$ws.Range("I136").Value = 71429540
$ws.Range("J136").Value = 2116
$ws.Range("K136").Value = 214288620
$ws.Range("L136").Value = 6348
$ws.Range("M136").Value = -214286070
$ws.Range("N136").Value = -11448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2364.5
$ws.Range("I86").Value = 2056.2778
$ws.Range("J86").Value = 3751.5
$ws.Range("K86").Value = 2056.2778
$ws.Range("L86").Value = 3751.5
$ws.Range("M86").Value = -933.2777999999998
$ws.Range("N86").Value = -5997.5

$ws.Range("H89").Value = 2364.5
$ws.Range("I89").Value = 2056.2778
$ws.Range("J89").Value = 3751.5
$ws.Range("K89").Value = 10281.389
$ws.Range("L89").Value = 18757.5
$ws.Range("M89").Value = -4665.388999999999
$ws.Range("N89").Value = -29989.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("L94").Value = 1000
$ws.Range("N94").Value = -1902

$ws.Range("H122").Value = 4914.92
$ws.Range("I122").Value = 5036.375
$ws.Range("K122").Value = 15109.125
$ws.Range("M122").Value = -12659.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1027.9354
$ws.Range("I5").Value = 1140.6666
$ws.Range("K5").Value = 3421.9998
$ws.Range("M5").Value = -3309.9998

$ws.Range("H12").Value = 216.3
$ws.Range("J12").Value = 159.14285
$ws.Range("L12").Value = 477.42855
$ws.Range("N12").Value = -823.4285500000001

$ws.Range("H39").Value = 3657.4285
$ws.Range("J39").Value = 3633.6667
$ws.Range("L39").Value = 10901.0001
$ws.Range("N39").Value = -11489.0001

$ws.Range("H92").Value = 229.08571
$ws.Range("I92").Value = 232.33333
$ws.Range("K92").Value = 696.99999
$ws.Range("M92").Value = 551.00001

$ws.Range("H122").Value = 1769.2
$ws.Range("J122").Value = 1899.0769
$ws.Range("L122").Value = 17091.6921
$ws.Range("N122").Value = -21991.6921

$ws.Range("H131").Value = 10417711
$ws.Range("I131").Value = 125000370
$ws.Range("J131").Value = 1105.875
$ws.Range("K131").Value = 375001110
$ws.Range("L131").Value = 3317.625
$ws.Range("M131").Value = -374996070
$ws.Range("N131").Value = -13397.625

$ws.Range("H135").Value = 1027.9354
$ws.Range("I135").Value = 1140.6666
$ws.Range("K135").Value = 10265.9994
$ws.Range("M135").Value = -7730.999400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17311196
$ws.Range("I70").Value = 17860692
$ws.Range("J70").Value = 16670118
$ws.Range("K70").Value = 17860692
$ws.Range("L70").Value = 16670118
$ws.Range("M70").Value = -17860422
$ws.Range("N70").Value = -16670658

$ws.Range("H73").Value = 17311196
$ws.Range("I73").Value = 17860692
$ws.Range("J73").Value = 16670118
$ws.Range("K73").Value = 17860692
$ws.Range("L73").Value = 16670118
$ws.Range("M73").Value = -17859756
$ws.Range("N73").Value = -16671990

$ws.Range("H135").Value = 35100
$ws.Range("J135").Value = 33444.445
$ws.Range("L135").Value = 33444.445
$ws.Range("N135").Value = -43584.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11063.75
$ws.Range("J40").Value = 11063.75
$ws.Range("L40").Value = 11063.75
$ws.Range("N40").Value = -11335.75

$ws.Range("H61").Value = 889.17645
$ws.Range("I61").Value = 861.06665
$ws.Range("J61").Value = 1100
$ws.Range("K61").Value = 861.06665
$ws.Range("L61").Value = 1100
$ws.Range("M61").Value = -659.06665
$ws.Range("N61").Value = -1504

$ws.Range("H113").Value = 889.17645
$ws.Range("I113").Value = 861.06665
$ws.Range("J113").Value = 1100
$ws.Range("K113").Value = 861.06665
$ws.Range("L113").Value = 1100
$ws.Range("M113").Value = 1308.93335
$ws.Range("N113").Value = -5440

$ws.Range("H132").Value = 2217.3125
$ws.Range("I132").Value = 2119.739
$ws.Range("J132").Value = 2466.6667
$ws.Range("K132").Value = 6359.217000000001
$ws.Range("L132").Value = 7400.000100000001
$ws.Range("M132").Value = -3829.217000000001
$ws.Range("N132").Value = -12460.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 13890469
$ws.Range("I122").Value = 19232574
$ws.Range("J122").Value = 997.8
$ws.Range("K122").Value = 57697722
$ws.Range("L122").Value = 2993.4
$ws.Range("M122").Value = -57695272
$ws.Range("N122").Value = -7893.4

$ws.Range("H132").Value = 2533.5217
$ws.Range("I132").Value = 2172.2632
$ws.Range("K132").Value = 6516.7896
$ws.Range("M132").Value = -3986.7896
